$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the "Price" (D) column cells being updated so that
# values such as "1.00", "0.480", "67.438.50" are stored as literal text
# (matching the inlineStr cells in the source workbook) instead of being
# auto-coerced into numbers/dates by Excel.
$dCells = @("D2","D3","D5","D6","D7","D9","D10","D11","D12","D13","D14","D15","D16","D17","D19","D20","D21","D23","D24","D25","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell.
$ws.Range("D2").Value = "67.438.50"
$ws.Range("E2").Value = "  -2.92%  "
$ws.Range("D3").Value = "3.500.27"
$ws.Range("E3").Value = "  -4.71%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "602.59"
$ws.Range("E5").Value = "  -2.99%  "
$ws.Range("D6").Value = "149.56"
$ws.Range("E6").Value = "  -5.98%  "
$ws.Range("D7").Value = "3.497.77"
$ws.Range("E7").Value = "  -4.70%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.480"
$ws.Range("E9").Value = "  -3.19%  "
$ws.Range("D10").Value = "0.138"
$ws.Range("E10").Value = "  -4.47%  "
$ws.Range("D11").Value = "6.95"
$ws.Range("E11").Value = "  -3.56%  "
$ws.Range("D12").Value = "0.421"
$ws.Range("E12").Value = "  -4.40%  "
$ws.Range("D13").Value = "0.0000218"
$ws.Range("E13").Value = "  -4.76%  "
$ws.Range("D14").Value = "4.093.64"
$ws.Range("E14").Value = "  -4.65%  "
$ws.Range("D15").Value = "31.53"
$ws.Range("E15").Value = "  -2.39%  "
$ws.Range("D16").Value = "3.502.24"
$ws.Range("E16").Value = "  -4.40%  "
$ws.Range("D17").Value = "67.339.03"
$ws.Range("E17").Value = "  -3.08%  "
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("D19").Value = "6.40"
$ws.Range("E19").Value = "  -1.56%  "
$ws.Range("D20").Value = "14.99"
$ws.Range("E20").Value = "  -5.77%  "
$ws.Range("D21").Value = "447.37"
$ws.Range("E21").Value = "  -4.58%  "
$ws.Range("E22").Value = "  -12.61%  "
$ws.Range("D23").Value = "0.619"
$ws.Range("E23").Value = "  -4.97%  "
$ws.Range("D24").Value = "77.38"
$ws.Range("E24").Value = "  -2.80%  "
$ws.Range("D25").Value = "0.0000129"
$ws.Range("E25").Value = "  +5.85%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "3.642.07"
$ws.Range("E27").Value = "  -4.68%  "
$ws.Range("D28").Value = "10.10"
$ws.Range("E28").Value = "  -9.27%  "
$ws.Range("D29").Value = "8.20"
$ws.Range("E29").Value = "  -5.32%  "
$ws.Range("D30").Value = "2.47"
$ws.Range("E30").Value = "  -5.41%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "1.54"
$ws.Range("E32").Value = "  -7.40%  "
$ws.Range("D33").Value = "0.164"
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("D34").Value = "25.69"
$ws.Range("E34").Value = "  -3.46%  "
$ws.Range("D35").Value = "6.08"
$ws.Range("E35").Value = "  -4.83%  "
$ws.Range("D36").Value = "3.491.80"
$ws.Range("E36").Value = "  -5.04%  "
$ws.Range("D37").Value = "1.83"
$ws.Range("E37").Value = "  -6.97%  "
$ws.Range("D38").Value = "7.97"
$ws.Range("E38").Value = "  -3.71%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "174.67"
$ws.Range("E41").Value = "  -2.26%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "2.19"
$ws.Range("E42").Value = "  -1.63%  "
$ws.Range("D43").Value = "0.0874"
$ws.Range("E43").Value = "  -2.07%  "
$ws.Range("D44").Value = "5.39"
$ws.Range("E44").Value = "  -6.98%  "
$ws.Range("D45").Value = "0.879"
$ws.Range("E45").Value = "  -4.84%  "
$ws.Range("D46").Value = "45.44"
$ws.Range("E46").Value = "  -2.65%  "
$ws.Range("D47").Value = "27.48"
$ws.Range("E47").Value = "  -6.19%  "
$ws.Range("D48").Value = "1.27"
$ws.Range("E48").Value = "  +6.35%  "
$ws.Range("D49").Value = "2.55"
$ws.Range("E49").Value = "  -5.53%  "
$ws.Range("D50").Value = "7.53"
$ws.Range("E50").Value = "  -4.06%  "
$ws.Range("D51").Value = "0.996"
$ws.Range("E51").Value = "  -4.04%  "
